$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet title to reflect new "through" date
$ws.Name = "Through 2022-11-20"
$ws.Range("B1").Value = "November 2022 (through November 20)"

# Update carjacking counts for November columns (B, M, X, AI, AT, BE, BP, CA)
# and one correction cell (D7) per the 2022-11-20 data refresh
$ws.Range("X2").Value = 5
$ws.Range("AT2").Value = 2
$ws.Range("BE2").Value = 4
$ws.Range("BP2").Value = 4
$ws.Range("M3").Value = 1
$ws.Range("X3").Value = 5
$ws.Range("B4").Value = 2
$ws.Range("M4").Value = 3
$ws.Range("B5").Value = 7
$ws.Range("M5").Value = 5
$ws.Range("X5").Value = 11
$ws.Range("AT5").Value = 5
$ws.Range("BE5").Value = 5
$ws.Range("BP5").Value = 4
$ws.Range("B6").Value = 4
$ws.Range("AT6").Value = 2
$ws.Range("BE6").Value = 2
$ws.Range("BP6").Value = 5
$ws.Range("B7").Value = 3
$ws.Range("D7").Value = 3
$ws.Range("M7").Value = 3
$ws.Range("X7").Value = 2
$ws.Range("AT7").Value = 1
$ws.Range("BE7").Value = 1
$ws.Range("CA7").Value = 1
$ws.Range("B8").Value = 2
$ws.Range("B9").Value = 2
$ws.Range("M9").Value = 5
$ws.Range("X9").Value = 2
$ws.Range("AT9").Value = 2
$ws.Range("B10").Value = 1
$ws.Range("B11").Value = 2
$ws.Range("X11").Value = 1
$ws.Range("M13").Value = 1
$ws.Range("AI14").Value = 2
$ws.Range("B15").Value = 3
$ws.Range("X15").Value = 4
$ws.Range("CA15").Value = 1
$ws.Range("M16").Value = 8
$ws.Range("X16").Value = 12
$ws.Range("AI16").Value = 2
$ws.Range("BE16").Value = 3
$ws.Range("B17").Value = 1
$ws.Range("X17").Value = 7
$ws.Range("BE17").Value = 2
$ws.Range("M18").Value = 3
$ws.Range("BE18").Value = 5
$ws.Range("M20").Value = 6
$ws.Range("BE20").Value = 1
$ws.Range("M21").Value = 8
$ws.Range("AI21").Value = 2
$ws.Range("X22").Value = 3
$ws.Range("BP23").Value = 4
$ws.Range("B24").Value = 3
$ws.Range("X24").Value = 6
$ws.Range("BP24").Value = 3
$ws.Range("B25").Value = 3
$ws.Range("X25").Value = 6
$ws.Range("AT25").Value = 5
$ws.Range("BE25").Value = 10
$ws.Range("BP25").Value = 3
$ws.Range("CA25").Value = 3
$ws.Range("M26").Value = 4
$ws.Range("M30").Value = 3
$ws.Range("B31").Value = 1
$ws.Range("X31").Value = 1
$ws.Range("M32").Value = 3
$ws.Range("AT32").Value = 1
$ws.Range("BE32").Value = 5
$ws.Range("B35").Value = 2
$ws.Range("B36").Value = 3
$ws.Range("X36").Value = 1
$ws.Range("B40").Value = 4
$ws.Range("M40").Value = 5
$ws.Range("X40").Value = 4
$ws.Range("AI40").Value = 3
$ws.Range("BP40").Value = 3
$ws.Range("B42").Value = 2
$ws.Range("BE42").Value = 1
$ws.Range("B44").Value = 1
$ws.Range("M44").Value = 1
$ws.Range("X44").Value = 3
$ws.Range("CA44").Value = 1
$ws.Range("X45").Value = 1
$ws.Range("M46").Value = 1
$ws.Range("X46").Value = 4
$ws.Range("B47").Value = 3
$ws.Range("X47").Value = 1
$ws.Range("AI47").Value = 1
$ws.Range("BE47").Value = 1
$ws.Range("AI48").Value = 1
$ws.Range("B50").Value = 1
$ws.Range("M50").Value = 3
$ws.Range("M53").Value = 3
$ws.Range("X53").Value = 1
$ws.Range("AT57").Value = 1
$ws.Range("X59").Value = 3
$ws.Range("BP59").Value = 1
$ws.Range("M61").Value = 1
$ws.Range("M62").Value = 1
$ws.Range("M64").Value = 2
$ws.Range("X64").Value = 2
$ws.Range("M66").Value = 1
$ws.Range("X72").Value = 1
$ws.Range("BE72").Value = 2
$ws.Range("M74").Value = 3
$ws.Range("BE75").Value = 2
$ws.Range("M76").Value = 5
$ws.Range("X76").Value = 4
$ws.Range("B81").Value = 2
$ws.Range("M81").Value = 3
$ws.Range("X81").Value = 1
$ws.Range("BE81").Value = 1
$ws.Range("M82").Value = 1
$ws.Range("BE84").Value = 1
$ws.Range("X85").Value = 1
$ws.Range("M89").Value = 1
$ws.Range("M90").Value = 1
$ws.Range("B91").Value = 1
$ws.Range("AT91").Value = 2
$ws.Range("AI95").Value = 1
$ws.Range("M97").Value = 1
$ws.Range("X97").Value = 2
$ws.Range("AT98").Value = 4
